$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'25.979.66"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').Value = "'1.740.10"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = "'247.67"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.47%  '
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = "'0.5023"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.81%  '
$ws.Range('D8').Value = "'0.2732"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('E9').Value = '  +1.54%  '
$ws.Range('E10').Value = '  +1.42%  '
$ws.Range('D11').Value = "'1.740.74"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('D12').Value = "'0.6540"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.00%  '
$ws.Range('D13').Value = "'15.17"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').Value = "'4.727"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.12%  '
$ws.Range('D15').Value = "'77.82"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').Value = "'25.996.85"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.61%  '
$ws.Range('D19').Value = "'11.88"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.53%  '
$ws.Range('D20').Value = "'0.000006856"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.23%  '
$ws.Range('D21').Value = "'4.614"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +9.02%  '
$ws.Range('D22').Value = "'1.963.04"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = "'8.757"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.47%  '
$ws.Range('E24').Value = '  +3.62%  '
$ws.Range('D25').Value = "'134.34"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.28%  '
$ws.Range('D26').Value = "'1.499"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('D27').Value = "'15.27"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.88%  '
$ws.Range('D28').Value = "'1.786"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.05%  '
$ws.Range('D29').Value = "'105.36"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  +2.41%  '
$ws.Range('D31').Value = "'0.08166"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('D32').Value = "'3.704"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.25%  '
$ws.Range('D33').Value = "'0.04742"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.19%  '
$ws.Range('D34').Value = "'2.667"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('D35').Value = "'0.9967"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.87%  '
$ws.Range('D36').Value = "'0.6138"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('D37').Value = "'2.758"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.75%  '
$ws.Range('D38').Value = "'0.01620"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.73%  '
$ws.Range('D39').Value = "'1.948"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.96%  '
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('D41').Value = "'100.87"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.28%  '
$ws.Range('D42').Value = "'0.8021"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.89%  '
$ws.Range('D43').Value = "'0.3917"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.33%  '
$ws.Range('D44').Value = "'5.017"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.74%  '
$ws.Range('E45').Value = '  +4.63%  '
$ws.Range('D46').Value = "'6.377"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.73%  '
$ws.Range('D47').Value = "'55.87"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.22%  '
$ws.Range('D48').Value = "'0.05294"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.58%  '
$ws.Range('D49').Value = "'30.94"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.86%  '
$ws.Range('D50').Value = "'7.661"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.23%  '
$ws.Range('D51').Value = "'0.3488"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.23%  '
